$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28; this shifts the existing rows 28-50 down to 29-51
# (and the sheet's used range grows from A1:R50 to A1:R51), matching the diff.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new record.
$ws.Range("A28").Value = 2
$ws.Range("B28").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 44587
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = 100112032
$ws.Range("G28").Value = "Zapallo italiano"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 240
$ws.Range("K28").Value = 10000
$ws.Range("L28").Value = 11000
$ws.Range("M28").Value = 10500
$ws.Range("N28").Value = '$/caja 60 unidades'
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 175
$ws.Range("Q28").Value = 60
$ws.Range("R28").Value = "Hortaliza"
